$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = 13.098
$ws.Range("E6").Value = 13.189
$ws.Range("D7").Value = -7.199
$ws.Range("A8").Value = -21.107
$ws.Range("E9").Value = 12.443
$ws.Range("A10").Value = -20.712
$ws.Range("E10").Value = 12.525
$ws.Range("A12").Value = -21.694
$ws.Range("B13").Value = 6.606
$ws.Range("A18").Value = -21.694
$ws.Range("D20").Value = -8.222
